$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: update E4:K4 values
$ws.Range("E4").Value = 7005
$ws.Range("F4").Value = 4674
$ws.Range("G4").Value = 4333
$ws.Range("H4").Value = 4191
$ws.Range("I4").Value = 3987
$ws.Range("J4").Value = 4166
$ws.Range("K4").Value = 4376

# Row 5: update E5:K5 values
$ws.Range("E5").Value = 2040
$ws.Range("F5").Value = 1984
$ws.Range("G5").Value = 1712
$ws.Range("H5").Value = 1572
$ws.Range("I5").Value = 1423
$ws.Range("J5").Value = 1682
$ws.Range("K5").Value = 2047

# Row 5 E5:K5 style change: copy style from F4 (which uses the border-less "right" style)
$ws.Range("F4").Copy()
$ws.Range("E5:K5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A3").Select()
